$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 45, pushing the existing rows 45-58 down to 47-60.
$ws.Range("A45:R46").EntireRow.Insert()

# Fill the common (unchanged) columns for the two brand-new rows, matching the
# rest of the data set (same market / region / category / variety / quality / origin / class).
$commonCols = @{
  A = 2
  B = "Comercializadora del Agro de Limarí"
  C = "Coquimbo"
  E = 4
  F = 100112030
  G = "Poroto granado"
  H = "Sin especificar"
  I = "Primera"
  O = "Provincia de Limarí"
  R = "Hortaliza"
}

foreach ($row in 45,46) {
  foreach ($col in $commonCols.Keys) {
    $ws.Cells.Item($row, [int][char]$col - [int][char]'A' + 1).Value2 = $commonCols[$col]
  }
}

# Row 45: new weekly entry
$ws.Range("D45").Value2 = 44559
$ws.Range("J45").Value2 = 520
$ws.Range("K45").Value2 = 9500
$ws.Range("L45").Value2 = 10000
$ws.Range("M45").Value2 = 9750
$ws.Range("N45").Value2 = "$/caja 15 kilos"
$ws.Range("P45").Value2 = 650
$ws.Range("Q45").Value2 = 15

# Row 46: new weekly entry
$ws.Range("D46").Value2 = 44559
$ws.Range("J46").Value2 = 560
$ws.Range("K46").Value2 = 16000
$ws.Range("L46").Value2 = 17000
$ws.Range("M46").Value2 = 16500
$ws.Range("N46").Value2 = "$/malla 25 kilos"
$ws.Range("P46").Value2 = 660
$ws.Range("Q46").Value2 = 25

# Keep the date column's display format consistent with the rest of column D.
$ws.Range("D45:D46").NumberFormat = $ws.Range("D47").NumberFormat

$ws.Range("A1").Select()
